$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table of data and
# both need the "想去人数" (F column) values updated for rows 2,3,4,5,7.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2214
    $ws.Range("F3").Value = 635
    $ws.Range("F4").Value = 1604
    $ws.Range("F5").Value = 7444
    $ws.Range("F7").Value = 194
}
